$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----
$ws.Cells.Item(1, 1).Value = "Tarefa"
$ws.Cells.Item(1, 2).Value = "Data"
$ws.Cells.Item(1, 3).Value = "Observação"

# ---- Data rows: Tarefa | Data | Observação ----
$data = @(
    @("NEONATURE",    "05/02/2024", "19 ETIQ"),
    @("GWS",           "05/02/2024", "525 ETIQ"),
    @("VANGUARDA",     "05/02/2024", "142 ETIQ"),
    @("ASTROMIC",      "05/02/2024", "196 ETIQ"),
    @("NEURO BETES",   "05/02/2024", "12 ETIQ"),
    @("HERA",          "05/02/2024", "19 ETIQ"),
    @("LUNO",          "05/02/2024", "385 ETIQ"),
    @("MF",            "05/02/2024", "SEM PEDIDOS"),
    @("RADT",          "05/02/2024", "SEM PEDIDOS"),
    @("DESAGITA",      "05/02/2024", "37 ETIQ"),
    @("DIABETINA",     "05/02/2024", "SEM PEDIDOS"),
    @("FITNESS",       "05/02/2024", "SEM PEDIDOS")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]

    # Column B holds a date-looking string ("05/02/2024") that must be stored
    # as literal text, not auto-converted to a date serial number. Force a
    # text format while assigning, then restore the cell to the default
    # "Normal" style so the saved cell carries no explicit style reference.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $item[1]
    $dateCell.Style = "Normal"

    $ws.Cells.Item($row, 3).Value = $item[2]

    $row++
}
